# feat: add 2022-Q4 data
#
# - "总计" sheet gets a new summary row for "2022-Q4" (inserted above the
#   existing "2022-Q3" summary row, which shifts down).
# - A new "2022-Q4" worksheet (fund holdings for the quarter) is inserted
#   between "总计" and "2022-Q3".

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)     # "总计"
$q3 = $wb.Worksheets.Item(2)        # "2022-Q3" (existing fund-holding sheet)

function Set-TextCell($rng, $val) {
    # Force a literal-text cell (keeps leading zeros / trailing zeros
    # like "016250" or "9.60" instead of Excel auto-coercing the
    # numeric-looking string to a real number), then drop back to the
    # default style so no stray number-format style index is left
    # behind on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift the existing summary row down and insert a new
#    2022-Q4 summary row above it.
# ---------------------------------------------------------------------

# Preserve formatting of row 2 (the "index" style on column A) onto the
# row it is about to become (row 3) before we overwrite row 2's values.
$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4122)

# Row 3 = old 2022-Q3 summary data (shifted down)
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.01

# Row 2 = new 2022-Q4 summary data
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.66

# ---------------------------------------------------------------------
# 2) Duplicate the existing "2022-Q3" sheet so the untouched fund-holding
#    data survives as the (new) "2022-Q3" sheet, then turn the original
#    sheet into "2022-Q4" and replace its contents with the Q4 data.
#    This keeps sheetId/r:id churn minimal: "2022-Q4" inherits the old
#    "2022-Q3" sheet's id/relationship, and the duplicate gets the new one.
# ---------------------------------------------------------------------

$q3.Copy($null, $q3)
$q3Dup = $wb.Worksheets.Item(3)
$q3.Name = "2022-Q4"
$q3Dup.Name = "2022-Q3"

$new = $q3
$new.Cells.Clear()

# Match the "总计" sheet's page margins (0.75/0.75/1/1/0.5/0.5 in, i.e.
# 54/54/72/72/36/36 pt) instead of the inherited "2022-Q3" margins.
$new.PageSetup.LeftMargin = 54
$new.PageSetup.RightMargin = 54
$new.PageSetup.TopMargin = 72
$new.PageSetup.BottomMargin = 72
$new.PageSetup.HeaderMargin = 36
$new.PageSetup.FooterMargin = 36

# Header row
Set-TextCell $new.Range("B1") "基金代码"
Set-TextCell $new.Range("C1") "基金名称"
Set-TextCell $new.Range("D1") "基金规模"
Set-TextCell $new.Range("E1") "股票总仓位"
Set-TextCell $new.Range("F1") "仓位占比"
Set-TextCell $new.Range("G1") "持有市值(亿元)"
Set-TextCell $new.Range("H1") "仓位排名"

# Row 2
$new.Range("A2").Value = 0
Set-TextCell $new.Range("B2") "016250"
Set-TextCell $new.Range("C2") "华夏远见成长一年持有混合A"
Set-TextCell $new.Range("D2") "9.60"
Set-TextCell $new.Range("E2") "88.62"
Set-TextCell $new.Range("F2") "3.77"
Set-TextCell $new.Range("G2") "0.3619"
$new.Range("H2").Value = 8

# Row 3
$new.Range("A3").Value = 1
Set-TextCell $new.Range("B3") "011404"
Set-TextCell $new.Range("C3") "融通鑫新成长混合C"
Set-TextCell $new.Range("D3") "5.57"
Set-TextCell $new.Range("E3") "94.42"
Set-TextCell $new.Range("F3") "2.36"
Set-TextCell $new.Range("G3") "0.1315"
$new.Range("H3").Value = 10

# Row 4
$new.Range("A4").Value = 2
Set-TextCell $new.Range("B4") "016251"
Set-TextCell $new.Range("C4") "华夏远见成长一年持有混合C"
Set-TextCell $new.Range("D4") "2.97"
Set-TextCell $new.Range("E4") "88.62"
Set-TextCell $new.Range("F4") "3.77"
Set-TextCell $new.Range("G4") "0.1120"
$new.Range("H4").Value = 8

# Row 5
$new.Range("A5").Value = 3
Set-TextCell $new.Range("B5") "011403"
Set-TextCell $new.Range("C5") "融通鑫新成长混合A"
Set-TextCell $new.Range("D5") "2.22"
Set-TextCell $new.Range("E5") "94.42"
Set-TextCell $new.Range("F5") "2.36"
Set-TextCell $new.Range("G5") "0.0524"
$new.Range("H5").Value = 10

# Apply the "header" / "index-column" style (matching the "总计" sheet's
# convention, same as the original workbook uses for its sheet headers).
$total.Range("B1:D1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$new.Range("A2:A5").PasteSpecial(-4122)

# Restore "2022-Q3" as the active sheet (unchanged from before the edit).
$q3Dup.Activate()
